# Se agrega parrafo 2
# Parrafo con formato despues de bandera
#
# Appends three new paragraphs after the existing "Araceli" paragraph:
#   1) an italic paragraph of filler text
#   2) a second italic paragraph of filler text
#   3) a trailing empty paragraph whose mark is bold
#
# We build the exact OOXML for the new paragraphs (so the rPr / pPr
# formatting matches precisely - rFonts Abadi, i/iCs for the text
# paragraphs, b/bCs for the trailing empty paragraph mark) and insert it
# at the end of the document body with Range.InsertXML, which splices raw
# WordprocessingML into the story without disturbing anything already
# there.

$d = $word.ActiveDocument

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$para1Text = "Lorem ipsum dolor sit amet, consectetur adipiscing elit. In facilisis diam mi, a gravida lacus congue sit amet. Ut fringilla vitae ipsum vel convallis. Nunc sed lacus dolor. Praesent tristique egestas ex, id faucibus tortor euismod ut. Duis sit amet dui tellus. Nam euismod justo urna, ut cursus mauris varius quis. Duis libero mi, condimentum quis volutpat non, gravida a lorem. Vivamus tincidunt dapibus neque, at porttitor arcu. Vivamus vel lacus et magna bibendum rutrum at consectetur odio. Praesent blandit eros a libero vehicula imperdiet. Etiam rhoncus suscipit aliquet. Sed ullamcorper leo quis lectus vestibulum imperdiet."
$para2Text = "In at mi ligula. Etiam vestibulum felis mauris, vitae sagittis lorem tincidunt a. Aliquam vitae ante dapibus, ullamcorper nisi sed, efficitur nisl. In eu sodales elit. Aliquam sed justo tortor. Vestibulum facilisis viverra nibh ut tincidunt. Etiam non odio in augue iaculis lacinia. Pellentesque elementum quis mi nec blandit. Phasellus sagittis, nulla at ullamcorper fermentum, dui lectus sodales lectus, at consectetur nulla lorem in felis. Cras nibh est, varius sed sapien quis, fringilla placerat dolor. In accumsan nulla quam, eu scelerisque ex consequat quis. Aliquam consequat et nunc vitae accumsan. Mauris rhoncus velit id lacinia suscipit. Donec pellentesque congue orci, in bibendum orci luctus nec. Suspendisse condimentum luctus dui, at tristique sem pellentesque a. Aenean venenatis, ligula ac bibendum sagittis, arcu mi condimentum leo, nec aliquam dolor leo id erat."

$newParasXml = @"
<w:p $wNs>
  <w:pPr>
    <w:rPr>
      <w:rFonts w:ascii="Abadi" w:hAnsi="Abadi"/>
      <w:i/>
      <w:iCs/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Abadi" w:hAnsi="Abadi"/>
      <w:i/>
      <w:iCs/>
    </w:rPr>
    <w:t>$para1Text</w:t>
  </w:r>
</w:p>
<w:p $wNs>
  <w:pPr>
    <w:rPr>
      <w:rFonts w:ascii="Abadi" w:hAnsi="Abadi"/>
      <w:i/>
      <w:iCs/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Abadi" w:hAnsi="Abadi"/>
      <w:i/>
      <w:iCs/>
    </w:rPr>
    <w:t>$para2Text</w:t>
  </w:r>
</w:p>
<w:p $wNs>
  <w:pPr>
    <w:rPr>
      <w:rFonts w:ascii="Abadi" w:hAnsi="Abadi"/>
      <w:b/>
      <w:bCs/>
    </w:rPr>
  </w:pPr>
</w:p>
"@

$end = $d.Content
$end.Collapse(0)
$result = $end.InsertXML($newParasXml)

Write-Output "Inserted new paragraphs; document now has $($d.Paragraphs.Count) paragraphs."
